$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("A3").Value = "Team Don't Panic"

$ws.Range("A4").Select()
